# Append: 2026-01-20 06:32 JST
# The scraper's rolling job-postings sheet ("ランサーズ") is refreshed:
# only 2 postings survive this run (rows 2-3), everything else (rows 4-15
# from the previous run) is dropped, and their per-row hyperlinks go with
# them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Drop the old rows (4-15) that didn't survive this refresh ---------
$ws.Range("A4:H15").EntireRow.Delete()

# Clear every hyperlink on the sheet; the two that remain relevant
# (F2, F3) get re-created below pointing at their (possibly new) URLs.
$ws.Cells.Hyperlinks.Delete()

# --- Row 2: replaced by a freshly scraped posting -----------------------
$ws.Range("A2").Value = "2026-01-20 06:32:45"
$ws.Range("B2").Value = "【急募】業務システムの開発・運用・保守エンジニア募集(フロントエンド/バックエンド)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5474899"
$ws.Range("G2").Value = 75
$ws.Range("H2").Value = "◆開発"

# --- Row 3: replaced by what used to be the last row of the old list ---
$ws.Range("A3").Value = "2026-01-20 06:32:45"
$ws.Range("B3").Value = "Google clab用マークシートCSV出力プログラム作成依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5474679"
$ws.Range("G3").Value = 10
$ws.Range("H3").ClearContents()

# --- Hyperlinks for the surviving URL cells ------------------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5474899")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5474679")

# --- Column width tweaks (B 52->44, D 30->26, H 13->12) ------------------
# ColumnWidth is in "characters"; Excel rounds it to whole pixels
# internally (pixels = round(width*6)) before re-deriving the stored
# OOXML <col width>, which is (pixels+5)/6. Biasing by -5/6 characters
# lands exactly back on the intended whole-character width after that
# round trip instead of drifting to N.8333.
$ws.Columns.Item(2).ColumnWidth = 44 - 5/6
$ws.Columns.Item(4).ColumnWidth = 26 - 5/6
$ws.Columns.Item(8).ColumnWidth = 12 - 5/6
